$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 50000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H99").Value = 796.4
$ws.Range("I99").Value = 538
$ws.Range("J99").Value = 1184
$ws.Range("K99").Value = 1614
$ws.Range("L99").Value = 3552
$ws.Range("M99").Value = -116
$ws.Range("N99").Value = -6548
$ws.Range("H103").Value = 357.1111
$ws.Range("J103").Value = 395
$ws.Range("L103").Value = 1185
$ws.Range("N103").Value = -2357
$ws.Range("H132").Value = 2865.121
$ws.Range("I132").Value = 2519.75
$ws.Range("K132").Value = 7559.25
$ws.Range("M132").Value = -5029.25
$ws.Range("H138").Value = 4365.886
$ws.Range("J138").Value = 4611.6553
$ws.Range("L138").Value = 13834.9659
$ws.Range("N138").Value = -24114.9659

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1299.1818
$ws.Range("I2").Value = 1532.5294
$ws.Range("K2").Value = 1532.5294
$ws.Range("M2").Value = -1419.5294
$ws.Range("H23").Value = 1100001
$ws.Range("J23").Value = 1100001
$ws.Range("L23").Value = 1100001
$ws.Range("N23").Value = -1100519
$ws.Range("H32").Value = 11130.278
$ws.Range("I32").Value = 4234.2744
$ws.Range("J32").Value = 46299.9
$ws.Range("K32").Value = 4234.2744
$ws.Range("L32").Value = 46299.9
$ws.Range("M32").Value = -3947.2744
$ws.Range("N32").Value = -46873.9
$ws.Range("H116").Value = 1299.1818
$ws.Range("I116").Value = 1532.5294
$ws.Range("K116").Value = 1532.5294
$ws.Range("M116").Value = 761.4706000000001
$ws.Range("H132").Value = 8665.415000000001
$ws.Range("I132").Value = 6454.963
$ws.Range("J132").Value = 19516.727
$ws.Range("K132").Value = 19364.889
$ws.Range("L132").Value = 58550.181
$ws.Range("M132").Value = -16834.889
$ws.Range("N132").Value = -63610.181

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1299.1818
$ws.Range("I3").Value = 1532.5294
$ws.Range("K3").Value = 1532.5294
$ws.Range("M3").Value = -1418.5294

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1875.08
$ws.Range("I122").Value = 1744.909
$ws.Range("K122").Value = 5234.727000000001
$ws.Range("M122").Value = -2784.727000000001
$ws.Range("H132").Value = 2600.0908
$ws.Range("I132").Value = 2109.96
$ws.Range("K132").Value = 6329.88
$ws.Range("M132").Value = -3799.88

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1161.6923
$ws.Range("I5").Value = 1067.25
$ws.Range("J5").Value = 1312.8
$ws.Range("K5").Value = 3201.75
$ws.Range("L5").Value = 3938.4
$ws.Range("M5").Value = -3089.75
$ws.Range("N5").Value = -4162.4
$ws.Range("H34").Value = 3288.3333
$ws.Range("J34").Value = 4892.75
$ws.Range("L34").Value = 14678.25
$ws.Range("N34").Value = -14846.25
$ws.Range("H39").Value = 7649.6665
$ws.Range("J39").Value = 7649.6665
$ws.Range("L39").Value = 22948.9995
$ws.Range("N39").Value = -23536.9995
$ws.Range("H68").Value = 1978.6666
$ws.Range("I68").Value = 1485
$ws.Range("J68").Value = 2077.4
$ws.Range("K68").Value = 4455
$ws.Range("L68").Value = 6232.200000000001
$ws.Range("M68").Value = -3644
$ws.Range("N68").Value = -7854.200000000001
$ws.Range("H71").Value = 1978.6666
$ws.Range("I71").Value = 1485
$ws.Range("J71").Value = 2077.4
$ws.Range("K71").Value = 13365
$ws.Range("L71").Value = 18696.6
$ws.Range("M71").Value = -9309
$ws.Range("N71").Value = -26808.6
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -7872
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 18000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -27360
$ws.Range("H86").Value = 929.6
$ws.Range("I86").Value = 966
$ws.Range("J86").Value = 875
$ws.Range("K86").Value = 2898
$ws.Range("L86").Value = 2625
$ws.Range("M86").Value = -1712
$ws.Range("N86").Value = -4997
$ws.Range("H87").Value = 15118.857
$ws.Range("I87").Value = 2502.8
$ws.Range("K87").Value = 7508.400000000001
$ws.Range("M87").Value = -6260.400000000001
$ws.Range("H89").Value = 929.6
$ws.Range("I89").Value = 966
$ws.Range("J89").Value = 875
$ws.Range("K89").Value = 8694
$ws.Range("L89").Value = 7875
$ws.Range("M89").Value = -2766
$ws.Range("N89").Value = -19731
$ws.Range("H90").Value = 15118.857
$ws.Range("I90").Value = 2502.8
$ws.Range("K90").Value = 22525.2
$ws.Range("M90").Value = -16285.2
$ws.Range("H117").Value = 108095
$ws.Range("I117").Value = 3779.6
$ws.Range("J117").Value = 212410.4
$ws.Range("K117").Value = 11338.8
$ws.Range("L117").Value = 637231.2
$ws.Range("M117").Value = -7896.799999999999
$ws.Range("N117").Value = -644115.2
$ws.Range("H122").Value = 4471
$ws.Range("J122").Value = 800
$ws.Range("L122").Value = 7200
$ws.Range("N122").Value = -12100
$ws.Range("H135").Value = 1161.6923
$ws.Range("I135").Value = 1067.25
$ws.Range("J135").Value = 1312.8
$ws.Range("K135").Value = 9605.25
$ws.Range("L135").Value = 11815.2
$ws.Range("M135").Value = -7070.25
$ws.Range("N135").Value = -16885.2
$ws.Range("H139").Value = 13380.951
$ws.Range("I139").Value = 16426.291
$ws.Range("K139").Value = 49278.87300000001
$ws.Range("M139").Value = -44138.87300000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19588.35
$ws.Range("I70").Value = 28113.348
$ws.Range("J70").Value = 8054.5293
$ws.Range("K70").Value = 28113.348
$ws.Range("L70").Value = 8054.5293
$ws.Range("M70").Value = -27843.348
$ws.Range("N70").Value = -8594.5293
$ws.Range("H73").Value = 19588.35
$ws.Range("I73").Value = 28113.348
$ws.Range("J73").Value = 8054.5293
$ws.Range("K73").Value = 28113.348
$ws.Range("L73").Value = 8054.5293
$ws.Range("M73").Value = -27177.348
$ws.Range("N73").Value = -9926.5293
$ws.Range("H122").Value = 4045.5
$ws.Range("I122").Value = 1943.2858
$ws.Range("K122").Value = 5829.857400000001
$ws.Range("M122").Value = -3379.857400000001
$ws.Range("H132").Value = 6922.3335
$ws.Range("I132").Value = 7213.3335
$ws.Range("J132").Value = 6194.8335
$ws.Range("K132").Value = 21640.0005
$ws.Range("L132").Value = 18584.5005
$ws.Range("M132").Value = -19110.0005
$ws.Range("N132").Value = -23644.5005

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 74998.5
$ws.Range("I5").Value = 74998.5
$ws.Range("K5").Value = 74998.5
$ws.Range("M5").Value = -74885.5
$ws.Range("H74").Value = 24500
$ws.Range("I74").Value = 24500
$ws.Range("K74").Value = 24500
$ws.Range("M74").Value = -23502
$ws.Range("H77").Value = 24500
$ws.Range("I77").Value = 24500
$ws.Range("K77").Value = 73500
$ws.Range("M77").Value = -68508
$ws.Range("H122").Value = 8065.815
$ws.Range("I122").Value = 5362.591
$ws.Range("K122").Value = 16087.773
$ws.Range("M122").Value = -13637.773
$ws.Range("H136").Value = 3990.86
$ws.Range("I136").Value = 5122.88
$ws.Range("K136").Value = 15368.64
$ws.Range("M136").Value = -12818.64

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1501.3112
$ws.Range("I132").Value = 1270.7693
$ws.Range("J132").Value = 2999.8333
$ws.Range("K132").Value = 3812.3079
$ws.Range("L132").Value = 8999.499899999999
$ws.Range("M132").Value = -1282.3079
$ws.Range("N132").Value = -14059.4999
$ws.Range("H136").Value = 5160.476
$ws.Range("I136").Value = 5229.9487
$ws.Range("K136").Value = 15689.8461
$ws.Range("M136").Value = -13139.8461
